$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "JUNK FOOD NIGHT"
$ws.Range("C2").Value = "write"
$ws.Range("D2").Value = "1000 food"
$ws.Range("E2").Value = "click"

# --- Row 3 ---
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "TROJAN HORSE FOR SALE"
$ws.Range("C3").Value = "write"
$ws.Range("D3").Value = " 1000 wood"
$ws.Range("E3").Value = "hotkey"

# --- Row 4 ---
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "ATM OF EREBUS"
$ws.Range("C4").Value = "write"
$ws.Range("D4").Value = "1000 gold"
$ws.Range("E4").Value = "write"

# --- Row 5 --- (E5 is already blank in the source sheet; leave it alone)
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "LAY OF THE LAND"
$ws.Range("C5").Value = "write"
$ws.Range("D5").Value = "Show map"

# --- Row 6 --- (E6 is already blank in the source sheet; leave it alone)
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "NINJACONNOR"
$ws.Range("C6").Value = "write"
$ws.Range("D6").Value = "100,000 each resource, maxes population cap, unlimited god power castings, 100x build/research speeds"

# --- Row 7 --- (E7 is already blank in the source sheet; leave it alone)
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "WRATH OF THE GODS"
$ws.Range("C7").Value = "write"
$ws.Range("D7").Value = "Get the Lightning Storm, Earthquake, Meteor and Tornado god powers"

# --- Row 9/10: make sure they exist as real (blank) rows, replacing the
# leftover style from the old B9 cell, matching the new blank layout ---
$ws.Range("A9:E9").Value = "__tmp__"
$ws.Range("A10:E10").Value = "__tmp__"
$ws.Range("A9:E10").ClearContents()
$ws.Range("A9:E10").Style = "Normal"

# --- Column widths ---
$ws.Columns.Item(2).ColumnWidth = 24.65
$ws.Columns.Item(3).ColumnWidth = 14.65
$ws.Columns.Item(4).ColumnWidth = 88

# --- View: zoom + selection ---
$win = $excel.ActiveWindow
$win.Zoom = 160
$ws.Range("D13").Select()
